# Insert two new weekly price records at rows 328-329 (Fruta / Plátano,
# Vega Modelo de Temuco), pushing the existing rows 328..377 down to
# 330..379.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("328:329").Insert()

# New row 328
$ws.Range("A328").Value = 10
$ws.Range("B328").Value = "Vega Modelo de Temuco"
$ws.Range("C328").Value = "La Araucanía"
$ws.Range("D328").Value = 44522
$ws.Range("E328").Value = 9
$ws.Range("F328").Value = "Fruta"
$ws.Range("G328").Value = 100108
$ws.Range("H328").Value = "Tropicales y subtropicales"
$ws.Range("I328").Value = 100108006
$ws.Range("J328").Value = "Plátano"
$ws.Range("K328").Value = "Sin especificar"
$ws.Range("L328").Value = "Maduro"
$ws.Range("M328").Value = 200
$ws.Range("N328").Value = 17000
$ws.Range("O328").Value = 17000
$ws.Range("P328").Value = 17000
$ws.Range("Q328").Value = "$/caja 20 kilos"
$ws.Range("R328").Value = "Ecuador"
$ws.Range("S328").Value = 850
$ws.Range("T328").Value = 20

# New row 329
$ws.Range("A329").Value = 10
$ws.Range("B329").Value = "Vega Modelo de Temuco"
$ws.Range("C329").Value = "La Araucanía"
$ws.Range("D329").Value = 44522
$ws.Range("E329").Value = 9
$ws.Range("F329").Value = "Fruta"
$ws.Range("G329").Value = 100108
$ws.Range("H329").Value = "Tropicales y subtropicales"
$ws.Range("I329").Value = 100108006
$ws.Range("J329").Value = "Plátano"
$ws.Range("K329").Value = "Sin especificar"
$ws.Range("L329").Value = "Pintón"
$ws.Range("M329").Value = 1200
$ws.Range("N329").Value = 21000
$ws.Range("O329").Value = 22000
$ws.Range("P329").Value = 21333
$ws.Range("Q329").Value = "$/caja 20 kilos"
$ws.Range("R329").Value = "Ecuador"
$ws.Range("S329").Value = 1067
$ws.Range("T329").Value = 20
